$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same style as other header cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-7: new columns I and J with numeric values
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 8
